# Reflect "new requirement separation": split the old single
# "Terms Typically Offered" column (D) into four columns:
#   D = Corequisites, E = Concurrent, F = Recommended, G = Terms Typically Offered
#
# Insert three blank columns at D:F. Excel shifts the existing D column
# (and its data/values) to G automatically, and the sheet's dimension is
# recalculated for us.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0

$ws.Columns("D:F").Insert()

# New header row
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# Default every data row's new Corequisites/Concurrent/Recommended cells to "NA"
$lastRow = 33
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Value = "NA"
    $ws.Cells.Item($r, 5).Value = "NA"
    $ws.Cells.Item($r, 6).Value = "NA"
}

# Row 24 (GSB 562) previously crammed "SP" (terms offered) and a
# "Corequisite: ..." note together into the old D column. Split them apart:
# the corequisite text becomes the real Corequisites (D) value, and the
# Terms Typically Offered (G) becomes just "SP".
$ws.Range("D24").Value = "OCOB graduate standing and GSB" + $nbsp + "511, GSB" + $nbsp + "513, GSB" + $nbsp + "523, GSB" + $nbsp + "531, GSB" + $nbsp + "533 and either GSB" + $nbsp + "512 or IME" + $nbsp + "503; and either GSB" + $nbsp + "524 or GSB" + $nbsp + "573; and either GSB" + $nbsp + "534 or IME" + $nbsp + "580, or approval from the Associate Dean."
$ws.Range("G24").Value = "SP"

# A handful of Prerequisites cells (column C) had their internal
# non-breaking spaces between a course prefix and number normalized to
# plain spaces as part of this pass.
$ws.Range("C10").Value = "GSB 512 or GSE 518 or IME 503."
$ws.Range("C16").Value = "GSA 555 or GSB 520."
$ws.Range("C21").Value = "GSE 520."
$ws.Range("C28").Value = "OCOB graduate standing; and GSB 512 or GSE 518 or IME 503 or approval from the Associate Dean."
$ws.Range("C31").Value = "GSB 531 and OCOB graduate standing or approval from the Associate Dean."
$ws.Range("C33").Value = "GSB 523 and OCOB graduate standing or approval from the Associate Dean."
